$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1) Reuse the existing blank-row formatting (styles 27/28/29) for the row
#    that becomes the new trailing blank row (row 34), by copying the
#    current row 30's format down onto it.
# ---------------------------------------------------------------------------
$ws.Range("A30:G30").Copy()
$ws.Range("A34:G34").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) Turn rows 30-33 into data rows matching the look of the other diary
#    entries (styles 19/20/20/20/20/20/21), by copying row 10's formatting.
# ---------------------------------------------------------------------------
$ws.Range("A10:G10").Copy()
$ws.Range("A30:G33").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3) Row heights for the new/changed rows.
# ---------------------------------------------------------------------------
$ws.Rows.Item(30).RowHeight = 227
$ws.Rows.Item(31).RowHeight = 136.85
$ws.Rows.Item(32).RowHeight = 214.5
$ws.Rows.Item(33).RowHeight = 72.2
$ws.Rows.Item(34).RowHeight = 15.5

# ---------------------------------------------------------------------------
# 4) Diary entry: Lecture 6 (new row 30) - 2/13/2020
# ---------------------------------------------------------------------------
$ws.Range("A30").Value = 43874
$ws.Range("B30").Value = "17:00 - 19:50"
$ws.Range("C30").Value = "None"
$ws.Range("D30").Value = "Finish exam on time, hope to remember what I learned. Hope we have the energy to concentrate on the lecture. "
$ws.Range("E30").Value = "Finished exam on time, remembered most of the details. Successfully sat through the lecture. The concepts were explained well and were easy to comprehend.  "
$ws.Range("F30").Value = "Felt like the time was just right, didn’t have time to go through the paper, should work at a faster pace for the finals. Handwriting got shabby towards the end, hopefully Kaj understands. The practical part of the exam was rather confusing to settle on. I think I spent a lot of time and couldn’t really write a satisfying answer. We discussed the survey results. Glad to know that a lot of people felt the same. Finally, we had a lecture about Stakeholders, developers and the bigger picture of the system. Really liked the KEP#9, because I tend to go for a run every time I am stressed and it has worked wonders. "
$ws.Range("G30").Value = "Tired. "

# ---------------------------------------------------------------------------
# 5) Diary entry: Team meeting (new row 31) - 2/17/2020
# ---------------------------------------------------------------------------
$ws.Range("A31").Value = 43878
$ws.Range("B31").Value = "16:00 - 18:00"
$ws.Range("C31").Value = "Team"
$ws.Range("D31").Value = "Discuss the functionalities of various features and developers"
$ws.Range("E31").Value = "We found a lot of features related to the functionalities and had to categorize it. We also managed to find the key developers in JabRef’s GitHub. "
$ws.Range("F31").Value = "After finding the functionalities, it was rather hard for each of us to agree whether it was functional or non-functional. For example, integration to environment, Built in custom export formats, etc. After agreeing on that, it was also hard for us to decide which four to put up for each. They all seemed equally crucial. "
$ws.Range("G31").Value = "Fruitful Discussion"

# ---------------------------------------------------------------------------
# 6) Diary entry: Team meeting (new row 32) - 2/19/2020
# ---------------------------------------------------------------------------
$ws.Range("A32").Value = 43880
$ws.Range("B32").Value = "21:00 -23:00"
$ws.Range("C32").Value = "Team"
$ws.Range("D32").Value = "Final discussion about functionalities, find stakeholder information, find what’s unique about the system."
$ws.Range("E32").Value = "We managed to find information about stakeholders after a tedious search of going among the donations, finding survey information related to JabRef. We were able to find unique features by reading reviews on JabRef "
$ws.Range("F32").Value = "After going through the survey details, we were somewhat able to settle on stakeholders. With a wide possible range of users mainly academics, we were able to find two universities that primarily used JabRef. We were happy that we could find this crucial info on time. After much thought, we finished the write up for features. While reading the reviews, we were leaning more towards what other people found useful about the system as well. This gave us a better sense of choosing the unique feature among our own contradicting beliefs. "
$ws.Range("G32").Value = "Happy that we could find related info which is very useful to write up the report. "

# ---------------------------------------------------------------------------
# 7) Diary entry: Add references (new row 33) - 2/20/2020
# ---------------------------------------------------------------------------
$ws.Range("A33").Value = 43881
$ws.Range("B33").Value = "7:30 - 9:00"
$ws.Range("C33").Value = "None"
$ws.Range("D33").Value = "Add references to the report"
$ws.Range("E33").Value = "Finished adding references to the report, made changes with grammar and was focussing on building a concise report."
$ws.Range("F33").Value = "Glad that we finished the report. Hopefully scores well. We have also included the references to ease the search for info for the grader. "
$ws.Range("G33").Value = "Satisfied!"
